$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that sits right after the
#    title paragraph at the top of the document.
# ---------------------------------------------------------------------
$metaRange = $d.Content
$metaRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metaRange.Expand(4) | Out-Null   # wdParagraph -> grow to the whole paragraph (incl. mark)
$metaRange.Delete()

# ---------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Alien Antix Free Online Slot
#    Review" right before the final (image-prompt) paragraph, by typing
#    the text + a paragraph break at the very start of that paragraph.
# ---------------------------------------------------------------------
$imgRange = $d.Content
$imgRange.Find.Execute('Create an eye-catching feature image', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lastParaStart = $imgRange.Start

$titleText = "Play Alien Antix Free Online Slot Review"
$insertRange = $d.Range($lastParaStart, $lastParaStart)
$insertRange.InsertAfter($titleText + "`r")

$newTitleRange = $d.Range($lastParaStart, $lastParaStart + $titleText.Length)
$newTitleRange.Font.Bold = 1

# ---------------------------------------------------------------------
# 3. Replace the old "Create an eye-catching..." image-prompt text with
#    the new meta-description text (keeping the existing italic run
#    formatting intact).
# ---------------------------------------------------------------------
$oldText = 'Create an eye-catching feature image for the online slot game "Alien Antix". The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. Make sure to incorporate elements of space and aliens in the image to match the theme of the game. The image should be vibrant and colorful, with the Maya warrior as the central focus, surrounded by aliens and other space objects. Make the image stand out to attract potential players to the game.'
$newText = "Experience unique gameplay with Alien Antix Slot, featuring bonuses and excellent graphics. Play for free and win big!"
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
